$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cxcl12"
$ws.Range("C2").Value = "Itga4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 109.026058
$ws.Range("H2").Value = 327.078174
$ws.Range("I2").Value = 0.3049840938689738
$ws.Range("J2").Value = 0.3049840938689738
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 26.04517333333333
$ws.Range("N2").Value = 78.13552
$ws.Range("O2").Value = 0.9210237118384171
$ws.Range("P2").Value = 0.921023711838417
$ws.Range("Q2").Value = 2839.602578460053
$ws.Range("R2").Value = 25556.42320614048
$ws.Range("S2").Value = 0.2808975821868785
$ws.Range("T2").Value = 0.2808975821868784

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cxcl12"
$ws.Range("C3").Value = "Itga4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 109.026058
$ws.Range("H3").Value = 327.078174
$ws.Range("I3").Value = 0.3049840938689738
$ws.Range("J3").Value = 0.3049840938689738
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.3302223333333333
$ws.Range("N3").Value = 0.990667
$ws.Range("O3").Value = 0.01167750336256582
$ws.Range("P3").Value = 0.01167750336256582
$ws.Range("Q3").Value = 36.00283926689533
$ws.Range("R3").Value = 324.025553402058
$ws.Range("S3").Value = 0.003561452781684032
$ws.Range("T3").Value = 0.003561452781684032

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cxcl12"
$ws.Range("C4").Value = "Itga4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 109.026058
$ws.Range("H4").Value = 327.078174
$ws.Range("I4").Value = 0.3049840938689738
$ws.Range("J4").Value = 0.3049840938689738
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.903109
$ws.Range("N4").Value = 5.709327
$ws.Range("O4").Value = 0.06729878479901708
$ws.Range("P4").Value = 0.06729878479901708
$ws.Range("Q4").Value = 207.488472214322
$ws.Range("R4").Value = 1867.396249928898
$ws.Range("S4").Value = 0.02052505890041129
$ws.Range("T4").Value = 0.02052505890041129

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cxcl12"
$ws.Range("C5").Value = "Itga4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 89.97721833333333
$ws.Range("H5").Value = 269.931655
$ws.Range("I5").Value = 0.2516978134001918
$ws.Range("J5").Value = 0.2516978134001917
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 26.04517333333333
$ws.Range("N5").Value = 78.13552
$ws.Range("O5").Value = 0.9210237118384171
$ws.Range("P5").Value = 0.921023711838417
$ws.Range("Q5").Value = 2343.472247542844
$ws.Range("R5").Value = 21091.2502278856
$ws.Range("S5").Value = 0.2318196543594579
$ws.Range("T5").Value = 0.2318196543594578

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Cxcl12"
$ws.Range("C6").Value = "Itga4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 89.97721833333333
$ws.Range("H6").Value = 269.931655
$ws.Range("I6").Value = 0.2516978134001918
$ws.Range("J6").Value = 0.2516978134001917
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.3302223333333333
$ws.Range("N6").Value = 0.990667
$ws.Range("O6").Value = 0.01167750336256582
$ws.Range("P6").Value = 0.01167750336256582
$ws.Range("Q6").Value = 29.71248698487611
$ws.Range("R6").Value = 267.4123828638849
$ws.Range("S6").Value = 0.002939202062331205
$ws.Range("T6").Value = 0.002939202062331204

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Cxcl12"
$ws.Range("C7").Value = "Itga4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 89.97721833333333
$ws.Range("H7").Value = 269.931655
$ws.Range("I7").Value = 0.2516978134001918
$ws.Range("J7").Value = 0.2516978134001917
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.903109
$ws.Range("N7").Value = 5.709327
$ws.Range("O7").Value = 0.06729878479901708
$ws.Range("P7").Value = 0.06729878479901708
$ws.Range("Q7").Value = 171.2364540051316
$ws.Range("R7").Value = 1541.128086046185
$ws.Range("S7").Value = 0.01693895697840266
$ws.Range("T7").Value = 0.01693895697840266

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Cxcl12"
$ws.Range("C8").Value = "Itga4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 158.477852
$ws.Range("H8").Value = 475.433556
$ws.Range("I8").Value = 0.4433180927308344
$ws.Range("J8").Value = 0.4433180927308344
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 26.04517333333333
$ws.Range("N8").Value = 78.13552
$ws.Range("O8").Value = 0.9210237118384171
$ws.Range("P8").Value = 0.921023711838417
$ws.Range("Q8").Value = 4127.583124834347
$ws.Range("R8").Value = 37148.24812350912
$ws.Range("S8").Value = 0.4083064752920808
$ws.Range("T8").Value = 0.4083064752920806

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Cxcl12"
$ws.Range("C9").Value = "Itga4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 158.477852
$ws.Range("H9").Value = 475.433556
$ws.Range("I9").Value = 0.4433180927308344
$ws.Range("J9").Value = 0.4433180927308344
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.3302223333333333
$ws.Range("N9").Value = 0.990667
$ws.Range("O9").Value = 0.01167750336256582
$ws.Range("P9").Value = 0.01167750336256582
$ws.Range("Q9").Value = 52.33292606909466
$ws.Range("R9").Value = 470.9963346218519
$ws.Range("S9").Value = 0.005176848518550587
$ws.Range("T9").Value = 0.005176848518550586

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Cxcl12"
$ws.Range("C10").Value = "Itga4"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 158.477852
$ws.Range("H10").Value = 475.433556
$ws.Range("I10").Value = 0.4433180927308344
$ws.Range("J10").Value = 0.4433180927308344
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.903109
$ws.Range("N10").Value = 5.709327
$ws.Range("O10").Value = 0.06729878479901708
$ws.Range("P10").Value = 0.06729878479901708
$ws.Range("Q10").Value = 301.600626441868
$ws.Range("R10").Value = 2714.405637976812
$ws.Range("S10").Value = 0.02983476892020312
$ws.Range("T10").Value = 0.02983476892020312
